$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values per regenerated save data (Strike# -> K)
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 0
